# Insert 3 new data rows for Mandarina "Murcott" at the top of the
# existing block of rows (rows 859-861), pushing the previous rows
# (and all rows after them) down by 3. This mirrors the commit's
# "new weekly data" addition pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 859; everything that
# was row 859 onward moves down to 862 onward.
$ws.Rows("859:861").Insert()

# Row 859: Murcott / Especial
$ws.Cells.Item(859, 1).Value = 9
$ws.Cells.Item(859, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(859, 3).Value = "Metropolitana"
$ws.Cells.Item(859, 4).Value = 45166
$ws.Cells.Item(859, 5).Value = 13
$ws.Cells.Item(859, 6).Value = "Fruta"
$ws.Cells.Item(859, 7).Value = 100102
$ws.Cells.Item(859, 8).Value = "Cítricos"
$ws.Cells.Item(859, 9).Value = 100102004
$ws.Cells.Item(859, 10).Value = "Mandarina"
$ws.Cells.Item(859, 11).Value = "Murcott"
$ws.Cells.Item(859, 12).Value = "Especial"
$ws.Cells.Item(859, 13).Value = 290
$ws.Cells.Item(859, 14).Value = 8000
$ws.Cells.Item(859, 15).Value = 8000
$ws.Cells.Item(859, 16).Value = 8000
$ws.Cells.Item(859, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(859, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(859, 19).Value = 800
$ws.Cells.Item(859, 20).Value = 10

# Row 860: Murcott / Primera
$ws.Cells.Item(860, 1).Value = 9
$ws.Cells.Item(860, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(860, 3).Value = "Metropolitana"
$ws.Cells.Item(860, 4).Value = 45166
$ws.Cells.Item(860, 5).Value = 13
$ws.Cells.Item(860, 6).Value = "Fruta"
$ws.Cells.Item(860, 7).Value = 100102
$ws.Cells.Item(860, 8).Value = "Cítricos"
$ws.Cells.Item(860, 9).Value = 100102004
$ws.Cells.Item(860, 10).Value = "Mandarina"
$ws.Cells.Item(860, 11).Value = "Murcott"
$ws.Cells.Item(860, 12).Value = "Primera"
$ws.Cells.Item(860, 13).Value = 250
$ws.Cells.Item(860, 14).Value = 6000
$ws.Cells.Item(860, 15).Value = 6000
$ws.Cells.Item(860, 16).Value = 6000
$ws.Cells.Item(860, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(860, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(860, 19).Value = 600
$ws.Cells.Item(860, 20).Value = 10

# Row 861: Murcott / Segunda
$ws.Cells.Item(861, 1).Value = 9
$ws.Cells.Item(861, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(861, 3).Value = "Metropolitana"
$ws.Cells.Item(861, 4).Value = 45166
$ws.Cells.Item(861, 5).Value = 13
$ws.Cells.Item(861, 6).Value = "Fruta"
$ws.Cells.Item(861, 7).Value = 100102
$ws.Cells.Item(861, 8).Value = "Cítricos"
$ws.Cells.Item(861, 9).Value = 100102004
$ws.Cells.Item(861, 10).Value = "Mandarina"
$ws.Cells.Item(861, 11).Value = "Murcott"
$ws.Cells.Item(861, 12).Value = "Segunda"
$ws.Cells.Item(861, 13).Value = 200
$ws.Cells.Item(861, 14).Value = 5000
$ws.Cells.Item(861, 15).Value = 5000
$ws.Cells.Item(861, 16).Value = 5000
$ws.Cells.Item(861, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(861, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(861, 19).Value = 500
$ws.Cells.Item(861, 20).Value = 10
